$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 96

# --- Column A: date, reuse the date style from the row above via PasteSpecial(formats) ---
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$ws.Cells.Item($row, 1).Value = 45597

# --- Columns B & C: plain numeric values ---
$ws.Cells.Item($row, 2).Value = 104.446907005456
$ws.Cells.Item($row, 3).Value = 121.06465164124

# --- Column D: text value "110.3" already exists elsewhere in the sheet (D93) ->
#     copy that cell verbatim so it lands as the very same shared string, no new style ---
$ws.Cells.Item(93, 4).Copy()
$ws.Cells.Item($row, 4).PasteSpecial(-4163)

# --- Column E: brand-new text value "112.3" -> stage it in a scratch cell forced to Text,
#     then copy its value+type into place (scratch cell is wiped afterwards) ---
$scratch = $ws.Cells.Item(1000, 1)
$scratch.NumberFormat = "@"
$scratch.Value = "112.3"
$scratch.Copy()
$ws.Cells.Item($row, 5).PasteSpecial(-4163)

# --- Column F: text value " 88.3" already exists elsewhere in the sheet (F94) ->
#     copy that cell verbatim so it reuses the existing shared string ---
$ws.Cells.Item(94, 6).Copy()
$ws.Cells.Item($row, 6).PasteSpecial(-4163)

# --- Column G: brand-new text value "169.7" -> same scratch-cell trick as column E ---
$scratch.NumberFormat = "@"
$scratch.Value = "169.7"
$scratch.Copy()
$ws.Cells.Item($row, 7).PasteSpecial(-4163)

# tidy up the scratch cell so it leaves no trace in the sheet
$scratch.Clear()
